$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.400.08"

$ws.Cells.Item(3, 4).Value = "2.650.08"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "598.01"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.01%  "

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "159.18"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +2.86%  "

$ws.Cells.Item(8, 5).Value = "  -1.02%  "

$ws.Cells.Item(9, 5).Value = "  -1.01%  "

$ws.Cells.Item(10, 5).Value = "  -1.14%  "

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "5.30"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +0.70%  "

$ws.Cells.Item(12, 5).Value = "  -0.89%  "

$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "28.02"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -0.28%  "

$ws.Cells.Item(14, 4).Value = "3.124.30"
$ws.Cells.Item(14, 5).Value = "  +0.03%  "

$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000188"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -2.50%  "

$ws.Cells.Item(16, 4).Value = "68.324.40"
$ws.Cells.Item(16, 5).Value = "  +0.22%  "

$ws.Cells.Item(17, 4).Value = "2.650.43"
$ws.Cells.Item(17, 5).Value = "  +0.42%  "

$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = "11.43"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +0.33%  "

$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "361.63"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -0.68%  "

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "7.44"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -0.69%  "

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "4.42"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.92%  "

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "4.77"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -2.51%  "

$ws.Cells.Item(23, 5).Value = "  +0.66%  "

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = "74.32"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.24%  "

$ws.Cells.Item(25, 5).Value = "  +0.06%  "

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "9.77"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  -0.42%  "

$ws.Cells.Item(27, 4).Value = "2.784.13"
$ws.Cells.Item(27, 5).Value = "  +0.26%  "

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = "0.0000104"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -2.41%  "

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  -0.01%  "

$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "562.26"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -1.74%  "

$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "8.04"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -1.17%  "

$ws.Cells.Item(32, 5).Value = "  -2.08%  "

$ws.Cells.Item(33, 5).Value = "  +0.78%  "

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "1.66"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +4.11%  "

$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.128"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.46%  "

$ws.Cells.Item(36, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.00%  "

$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "160.54"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -0.32%  "

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "19.67"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +1.52%  "

$ws.Cells.Item(39, 5).Value = "  -1.24%  "

$ws.Cells.Item(40, 5).Value = "  -1.38%  "

$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "5.33"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  -0.83%  "

$ws.Cells.Item(42, 5).Value = "  -1.87%  "

$ws.Cells.Item(43, 4).Value = "0.0₆0322"

$ws.Cells.Item(44, 5).Value = "  +0.01%  "

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "158.19"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +1.39%  "

$ws.Cells.Item(46, 5).Value = "  +0.98%  "

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "22.05"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +1.14%  "

$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.70"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -1.19%  "

$ws.Cells.Item(49, 5).Value = "  -2.15%  "

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.575"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +1.33%  "

$ws.Cells.Item(51, 5).Value = "  -0.99%  "
